# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps for files that were
# re-handed-off, and rolls those new timestamps up into the Overview
# sheet's "Latest Handoff Date" column.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# New handoff timestamps produced by this run.
$zhTimestamp = "2016-03-23 14:30:10"
$deTimestamp = "2016-03-23 14:30:17"

# zh-cn: column E = "Latest Handoff Datetime" for rows 7, 10-16
$zh.Range("E7").Value = $zhTimestamp
$zh.Range("E10").Value = $zhTimestamp
$zh.Range("E11").Value = $zhTimestamp
$zh.Range("E12").Value = $zhTimestamp
$zh.Range("E13").Value = $zhTimestamp
$zh.Range("E14").Value = $zhTimestamp
$zh.Range("E15").Value = $zhTimestamp
$zh.Range("E16").Value = $zhTimestamp

# de-de: column E = "Latest Handoff Datetime" for rows 7, 10-16
$de.Range("E7").Value = $deTimestamp
$de.Range("E10").Value = $deTimestamp
$de.Range("E11").Value = $deTimestamp
$de.Range("E12").Value = $deTimestamp
$de.Range("E13").Value = $deTimestamp
$de.Range("E14").Value = $deTimestamp
$de.Range("E15").Value = $deTimestamp
$de.Range("E16").Value = $deTimestamp

# Overview: column D = "Latest Handoff Date" (max across languages) for
# the two rows whose source files picked up a new handoff in this run.
$overview.Range("D14").Value = $deTimestamp
$overview.Range("D15").Value = $deTimestamp
